$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend header row (row 1) into new columns P and Q, copying the
#     existing header style (bold font + border + centered) from O1 ---
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Cells.Item(1, 16).Value = 14
$ws.Cells.Item(1, 17).Value = 15

# --- Fill new P/Q columns (rows 2-25) with value 2 (no special style,
#     matching the plain body cells) ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 16).Value = 2
    $ws.Cells.Item($r, 17).Value = 2
}

# --- Update the I/K/M/O body columns (rows 2-25): swap the 1/2 pattern
#     I: 1 -> 2 ; K: 2 -> 1 ; M: 1 -> 2 ; O: 2 -> 1 ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2
    $ws.Cells.Item($r, 11).Value = 1
    $ws.Cells.Item($r, 13).Value = 2
    $ws.Cells.Item($r, 15).Value = 1
}
